$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation_log_key")

# Insert two new rows at row 20, pushing the existing rows (old 20-22, 33-35) down.
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# Set column A values first (new shared strings 197, 198) so they land before
# the descriptive text strings (199, 200) in the shared-string table.
$ws.Range("A20").Value = "scenario_3.4.1"
$ws.Range("A21").Value = "scenario_3.4.2"

# Populate new row 20 (scenario_3.4.1)
$ws.Range("C20").Value = "Biennial breeding; Liz's biennial model; psi = 0.75, lambda in model w/ tight prior"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = "biennial"
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "Uniform(0.95, 1.05)"
$ws.Range("I20").Value = "Uniform 0.5-0.95"
$ws.Range("J20").Value = "Liz's derivation"
$ws.Range("K20").Value = "all three"

# Populate new row 21 (scenario_3.4.2)
$ws.Range("C21").Value = "Biennial breeding; Ben's biennial model; psi = 0.75, lambda in model w/ tight prior"
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = "biennial"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = "Uniform(0.95, 1.05)"
$ws.Range("I21").Value = "Uniform 0.5-0.95"
$ws.Range("J21").Value = "Liz's derivation"
$ws.Range("K21").Value = "all three"

# Remove the fill/highlight for column A on the two new rows ("No Fill")
$ws.Range("A20:A21").Interior.ColorIndex = -4142

# Update sheet view (scrolled position / selection) to match the saved state
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D21:K21").Select()

# Update workbook window size/position
$excel.ActiveWindow.WindowState = -4143
$excel.Windows.Item(1).Left = -120
$excel.Windows.Item(1).Top = -120
$excel.Windows.Item(1).Width = 20730
$excel.Windows.Item(1).Height = 11160
